$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.467.79"
$ws.Range("E2").Value = "  -2.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.224.66"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.23"
$ws.Range("E5").Value = "  -7.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.31"
$ws.Range("E6").Value = "  +11.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -3.92%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.08"
$ws.Range("E10").Value = "  -7.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  -2.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.22"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.83"
$ws.Range("E13").Value = "  -3.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "  +11.70%  "

$ws.Range("E15").Value = "  -2.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.13"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.556.64"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.222.47"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.397.98"
$ws.Range("E19").Value = "  -2.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  +7.65%  "

$ws.Range("E21").Value = "  -3.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.38"
$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  +20.20%  "

$ws.Range("E24").Value = "  -3.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "228.97"
$ws.Range("E25").Value = "  -2.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.15"
$ws.Range("E26").Value = "  -3.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.69"
$ws.Range("E27").Value = "  -2.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.34"
$ws.Range("E30").Value = "  -8.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.34"
$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("E32").Value = "  -5.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.02"
$ws.Range("E33").Value = "  -2.60%  "

$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.63"
$ws.Range("E35").Value = "  -1.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +11.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.37"
$ws.Range("E37").Value = "  +4.10%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0385"
$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.126"
$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.43"
$ws.Range("E41").Value = "  -4.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.10"
$ws.Range("E42").Value = "  -2.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.234"
$ws.Range("E43").Value = "  -1.03%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.62"
$ws.Range("E45").Value = "  -8.72%  "

$ws.Range("E46").Value = "  -4.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.43"
$ws.Range("E47").Value = "  -6.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("E48").Value = "  +3.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.50"
$ws.Range("E49").Value = "  +1.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.44"
$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("E51").Value = "  +6.00%  "

